$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controllers")
$r = $ws.Range("G10")
$r.Validation.Add(3, 1, 1, "=1,2,3")
$val = $r.Validation
$props = $val | Get-Member
foreach ($p in $props) { "$p" }
